$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure the Price column keeps its original text representation instead of
# being auto-coerced into a floating point number by Excel.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "67.234.06"
$ws.Range("E2").Value = "  -0.12%  "
$ws.Range("D3").Value = "2.542.68"
$ws.Range("E3").Value = "  -3.01%  "
$ws.Range("E4").Value = "  -0.03%  "
$ws.Range("D5").Value = "588.35"
$ws.Range("E5").Value = "  -0.62%  "
$ws.Range("D6").Value = "173.32"
$ws.Range("E6").Value = "  +3.82%  "
$ws.Range("E7").Value = "  +0.00%  "
$ws.Range("D8").Value = "0.528"
$ws.Range("E8").Value = "  -0.93%  "
$ws.Range("D9").Value = "2.540.28"
$ws.Range("E9").Value = "  -3.07%  "
$ws.Range("D10").Value = "0.140"
$ws.Range("E10").Value = "  +1.18%  "
$ws.Range("E11").Value = "  +0.34%  "
$ws.Range("D12").Value = "0.352"
$ws.Range("E12").Value = "  -3.64%  "
$ws.Range("E13").Value = "  -1.05%  "
$ws.Range("D14").Value = "27.04"
$ws.Range("E14").Value = "  -1.35%  "
$ws.Range("B15").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C15").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D15").Value = "2.999.22"
$ws.Range("E15").Value = "  -3.27%  "
$ws.Range("B16").Value = "ShibaInu"
$ws.Range("C16").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D16").Value = "0.0000178"
$ws.Range("E16").Value = "  -1.39%  "
$ws.Range("D17").Value = "67.066.96"
$ws.Range("E17").Value = "  -0.39%  "
$ws.Range("D18").Value = "2.541.01"
$ws.Range("E18").Value = "  -2.75%  "
$ws.Range("D19").Value = "8.09"
$ws.Range("E19").Value = "  +2.26%  "
$ws.Range("D20").Value = "11.42"
$ws.Range("E20").Value = "  -3.73%  "
$ws.Range("D21").Value = "354.16"
$ws.Range("E21").Value = "  -0.43%  "
$ws.Range("D22").Value = "4.23"
$ws.Range("E22").Value = "  -1.70%  "
$ws.Range("D23").Value = "4.67"
$ws.Range("E23").Value = "  +0.38%  "
$ws.Range("D24").Value = "1.98"
$ws.Range("E24").Value = "  +2.45%  "
$ws.Range("E25").Value = "  -0.03%  "
$ws.Range("D26").Value = "70.28"
$ws.Range("E26").Value = "  +1.65%  "
$ws.Range("D27").Value = "9.88"
$ws.Range("E27").Value = "  -4.97%  "
$ws.Range("D28").Value = "1.00"
$ws.Range("E28").Value = "  -0.37%  "
$ws.Range("D29").Value = "2.668.15"
$ws.Range("E29").Value = "  -3.17%  "
$ws.Range("D30").Value = "0.0₃0993"
$ws.Range("E30").Value = "  -1.24%  "
$ws.Range("B31").Value = "InternetComputer(DFINITY)"
$ws.Range("C31").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D31").Value = "8.32"
$ws.Range("E31").Value = "  +4.77%  "
$ws.Range("B32").Value = "Bittensor"
$ws.Range("C32").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D32").Value = "532.08"
$ws.Range("E32").Value = "  -2.15%  "
$ws.Range("E33").Value = "  -1.02%  "
$ws.Range("E34").Value = "  -1.41%  "
$ws.Range("E35").Value = "  -2.56%  "
$ws.Range("D36").Value = "0.999"
$ws.Range("E36").Value = "  -0.04%  "
$ws.Range("D37").Value = "1.47"
$ws.Range("E37").Value = "  -1.83%  "
$ws.Range("D38").Value = "157.79"
$ws.Range("E38").Value = "  +0.46%  "
$ws.Range("D39").Value = "18.75"
$ws.Range("E39").Value = "  -1.02%  "
$ws.Range("D40").Value = "18.44"
$ws.Range("E40").Value = "  +1.10%  "
$ws.Range("D41").Value = "0.357"
$ws.Range("E41").Value = "  -2.14%  "
$ws.Range("D42").Value = "1.81"
$ws.Range("E42").Value = "  +0.11%  "
$ws.Range("D43").Value = "5.16"
$ws.Range("E43").Value = "  -0.28%  "
$ws.Range("D44").Value = "2.52"
$ws.Range("E44").Value = "  +4.02%  "
$ws.Range("E45").Value = "  -0.03%  "
$ws.Range("B46").Value = "OKB"
$ws.Range("C46").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D46").Value = "39.76"
$ws.Range("E46").Value = "  -0.91%  "
$ws.Range("B47").Value = "Aave"
$ws.Range("C47").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D47").Value = "151.40"
$ws.Range("E47").Value = "  -0.40%  "
$ws.Range("B48").Value = "ARBITRUM"
$ws.Range("C48").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D48").Value = "0.564"
$ws.Range("E48").Value = "  -2.60%  "
$ws.Range("B49").Value = "Filecoin"
$ws.Range("C49").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D49").Value = "3.72"
$ws.Range("E49").Value = "  -1.81%  "
$ws.Range("D50").Value = "0.0₆0276"
$ws.Range("E50").Value = "  -9.23%  "
$ws.Range("B51").Value = "Optimism"
$ws.Range("C51").Value = "https://coinranking.com/coin/n1p-s_gm1+optimism-op"
$ws.Range("D51").Value = "1.73"
$ws.Range("E51").Value = "  +1.41%  "
